# Applies the "Add description to all cards" edit:
#  - Fills in the Manacost (column D) for rows 67-116 which were previously empty
#  - Fixes the Manacost value in D48 (6 -> 8)
#  - Updates the Type (column C) for 4 spell cards
#  - Rewrites 9 card description texts (column I) with updated/fixed wording
#  - Restores the last-known selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Manacost (column D) fix for row 48 ---------------------------------
$ws.Cells.Item(48, 4).Value = 8

# --- 2. Manacost (column D) values newly entered for rows 67-116 -----------
$manacost = @{
    67  = 5
    68  = 2
    69  = 4
    70  = 1
    71  = 2
    72  = 4
    73  = 7
    74  = 2
    75  = 0
    76  = 4
    77  = 1
    78  = 2
    79  = 3
    80  = 3
    81  = 7
    82  = 9
    83  = 3
    84  = 5
    85  = 2
    86  = 6
    87  = 1
    88  = 2
    89  = 6
    90  = 5
    91  = 1
    92  = 0
    93  = 3
    94  = 3
    95  = 2
    96  = 5
    97  = 3
    98  = 1
    99  = 6
    100 = 3
    101 = 6
    102 = 3
    103 = 7
    104 = 10
    105 = 5
    106 = 4
    107 = 10
    108 = 1
    109 = 2
    110 = 5
    111 = 2
    112 = 8
    113 = 1
    114 = 5
    115 = 2
    116 = 1
}

foreach ($row in $manacost.Keys) {
    $ws.Cells.Item($row, 4).Value = $manacost[$row]
}

# --- 3. Type (column C) corrections -----------------------------------------
$typeFixes = @{
    88  = "Spell master"
    99  = "Spell master"
    101 = "Spell hero"
    109 = "Spell hero"
}

foreach ($row in $typeFixes.Keys) {
    $ws.Cells.Item($row, 3).Value = $typeFixes[$row]
}

# --- 4. Description (column I) text rewrites --------------------------------
$descriptionFixes = @{
    71  = "Deal 4 damage to two random enemy minions"
    72  = "Deal 5 damage to a minion and the enemy"
    76  = "Deals 3-10 damage"
    90  = "Deal 10 damage to the minion"
    93  = "Shoot three missiles at random enemies that deal 4 damage each."
    98  = "Deal 2 damage to a minion and give it +3 attack"
    106 = "Choose a minion. Whenever it attacks, restore 4 health to your hero"
    116 = "Force your opponent to change a spell card with you."
    117 = "Open 3 cards from your opponent's hand"
}

foreach ($row in $descriptionFixes.Keys) {
    $ws.Cells.Item($row, 9).Value = $descriptionFixes[$row]
}

# --- 5. Restore the scroll position / selection ------------------------------
$excel.ActiveWindow.ScrollRow = 76
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("H114").Select()
